{"js": "// Insert the contact-info line as its own centered paragraph, directly\n// after the \"Dheeraj Chand\" title paragraph (and before the\n// \"PROFESSIONAL SUMMARY\" heading), matching the target diff exactly:\n//   <w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>\n//     <w:r><w:t>202.550.7110 | ... | Austin, TX</w:t></w:r>\n//   </w:p>\n//\n// We use insertOoxml (rather than insertParagraph + formatting tweaks)\n// because insertParagraph copies the anchor run/paragraph's direct\n// formatting (bold + 28pt from the \"Dheeraj Chand\" title run, or the\n// Heading2 paragraph style from \"PROFESSIONAL SUMMARY\"), which this\n// new paragraph must NOT have \u2014 it needs a bare run with no rPr and a\n// pPr containing only center justification.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0]; // \"Dheeraj Chand\"\nconst insertionRange = titleParagraph.getRange(Word.RangeLocation.after);\n\nconst contactText =\n  \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\";\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData></pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr><w:r><w:t>' +\n  contactText +\n  '</w:t></w:r></w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part>' +\n  '</pkg:package>';\n\ninsertionRange.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert the contact-info line as its own centered paragraph, directly\n# after the \"Dheeraj Chand\" title paragraph (and before the\n# \"PROFESSIONAL SUMMARY\" heading), matching the target diff exactly:\n#   <w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>\n#     <w:r><w:t>202.550.7110 | ... | Austin, TX</w:t></w:r>\n#   </w:p>\n#\n# Plain InsertParagraphAfter()/new-paragraph-then-set-Text approaches end\n# up copying the title run's direct formatting (bold + 28pt) and/or an\n# inherited paragraph style onto the new paragraph, none of which the\n# target has. Instead we replace the title paragraph's own Range (which\n# already spans through its trailing paragraph mark) with OOXML that\n# reproduces that same title paragraph unchanged followed by the new\n# contact paragraph - InsertXML() replaces exactly the Range it is\n# called on and nothing beyond it, so the next paragraph\n# (\"PROFESSIONAL SUMMARY\") is left completely untouched.\n\n$d = $word.ActiveDocument\n\n$titleParagraph = $d.Paragraphs.Item(1)   # \"Dheeraj Chand\"\n$titleRange = $d.Range($titleParagraph.Range.Start, $titleParagraph.Range.End)\n\n$contactText = \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\"\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/></w:rPr><w:t>Dheeraj Chand</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr><w:r><w:t>' + $contactText + '</w:t></w:r></w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part>' +\n  '</pkg:package>'\n\n$titleRange.InsertXML($ooxml)\n"}
